$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully clear column E (content + formatting), which is not part of the new table
$ws.Range("E1:E21").Clear()

# Clear contents only for A2:D21 (keep default formatting, matches the target)
$ws.Range("A2:D21").ClearContents()

# Header row (A1:D1 already carry the bold/border header style - only update the text)
$ws.Range("A1").Value = "Pergunta"
$ws.Range("B1").Value = "Onix"
$ws.Range("C1").Value = "HB20S"
$ws.Range("D1").Value = "Virtus"

# Data rows 2-21
$ws.Range("A2").Value = "Alerta de Ponto Cego "
$ws.Range("B2").Value = "x"
$ws.Range("C2").Value = "x"

$ws.Range("A3").Value = "Sensores de estacionamento (dianteiro, lateral e traseiro) "
$ws.Range("B3").Value = "x"

$ws.Range("A4").Value = " Alerta de frenagem de emergência (ESS)"
$ws.Range("B4").Value = "x"
$ws.Range("C4").Value = "x"

$ws.Range("A5").Value = "Sistema de frenagem automática pós colisão `"Post Collision Brake`". "
$ws.Range("D5").Value = "x"

$ws.Range("A6").Value = " Alerta de colisão frontal e frenagem automática "
$ws.Range("C6").Value = "x"
$ws.Range("D6").Value = "x"

$ws.Range("A7").Value = "Assistente de permanência em faixa"
$ws.Range("C7").Value = "x"

$ws.Range("A8").Value = "Farol alto adaptativo "
$ws.Range("C8").Value = "x"

$ws.Range("A9").Value = "Alerta de tráfego cruzado traseiro"
$ws.Range("C9").Value = "x"

$ws.Range("A10").Value = "6 airbags"
$ws.Range("B10").Value = "x"
$ws.Range("C10").Value = "x"
$ws.Range("D10").Value = "x"

$ws.Range("A11").Value = "Sensor crepuscular"
$ws.Range("B11").Value = "x"
$ws.Range("C11").Value = "x"
$ws.Range("D11").Value = "x"

$ws.Range("A12").Value = "Alerta de abertura de porta (Alerta de saída segura)"
$ws.Range("C12").Value = "x"

$ws.Range("A13").Value = "Easy Park - Assistente de estacionamento automático"
$ws.Range("B13").Value = "x"

$ws.Range("A14").Value = "Retrovisores externos com rebatimento elétrico"
$ws.Range("C14").Value = "x"
$ws.Range("D14").Value = "x"

$ws.Range("A15").Value = "Paddle Shift"
$ws.Range("C15").Value = "x"
$ws.Range("D15").Value = "x"

$ws.Range("A16").Value = "Sistema Stop/Start"
$ws.Range("C16").Value = "x"
$ws.Range("D16").Value = "x"

$ws.Range("A17").Value = "Multimidia de 11”"
$ws.Range("B17").Value = "x"

$ws.Range("A18").Value = "OnStar"
$ws.Range("B18").Value = "x"

$ws.Range("A19").Value = "Wi-Fi embarcado"
$ws.Range("B19").Value = "x"

$ws.Range("A20").Value = "Projeção (Android Auto e Apple CarPlay) sem o uso de cabo"
$ws.Range("B20").Value = "x"
$ws.Range("C20").Value = "x"

$ws.Range("A21").Value = "Painel de instrumentos digital de 10.25"
$ws.Range("D21").Value = "x"

# Update the selection to match the final state
$ws.Range("A17:D21").Select()
